$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 596 (existing data for rows 596:607 shifts down to 599:610)
$ws.Rows.Item(596).EntireRow.Insert()
$ws.Rows.Item(596).EntireRow.Insert()
$ws.Rows.Item(596).EntireRow.Insert()

# Populate the 3 new rows (596:598) with the new week's data (week of 44656),
# mirroring the structure of the surrounding rows.
$rows = @(596, 597, 598)
$labels = @("Pintón", "Primera Maduro", "Primera Pintón")
$vols = @(120, 160, 160)
$prices = @(16000, 18000, 19000)
$kgPrices = @(800, 900, 950)

for ($i = 0; $i -lt 3; $i++) {
    $r = $rows[$i]

    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44656
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108006
    $ws.Cells.Item($r, 10).Value = "Plátano"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $labels[$i]
    $ws.Cells.Item($r, 13).Value = $vols[$i]
    $ws.Cells.Item($r, 14).Value = $prices[$i]
    $ws.Cells.Item($r, 15).Value = $prices[$i]
    $ws.Cells.Item($r, 16).Value = $prices[$i]
    $ws.Cells.Item($r, 17).Value = "$/caja 20 kilos"
    $ws.Cells.Item($r, 18).Value = "Ecuador"
    $ws.Cells.Item($r, 19).Value = $kgPrices[$i]
    $ws.Cells.Item($r, 20).Value = 20
}
